# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Mon Dec 11 17:47:51 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.626.74"
$ws.Range("E2").Value = "  -4.92%  "

$ws.Range("D3").Value = "2.209.35"
$ws.Range("E3").Value = "  -5.66%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").Value = "244.23"
$ws.Range("E5").Value = "  +2.25%  "

$ws.Range("E6").Value = "  -5.69%  "

$ws.Range("D7").Value = "70.41"
$ws.Range("E7").Value = "  -3.01%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("E9").Value = "  -8.04%  "

$ws.Range("D10").Value = "36.97"
$ws.Range("E10").Value = "  +11.75%  "

$ws.Range("D11").Value = "0.0948"
$ws.Range("E11").Value = "  -5.46%  "

$ws.Range("D12").Value = "57.78"
$ws.Range("E12").Value = "  -4.72%  "

$ws.Range("E13").Value = "  -2.87%  "

$ws.Range("E14").Value = "  -7.27%  "

$ws.Range("D15").Value = "2.537.79"
$ws.Range("E15").Value = "  -5.79%  "

$ws.Range("D16").Value = "14.72"
$ws.Range("E16").Value = "  -8.16%  "

$ws.Range("D17").Value = "0.839"
$ws.Range("E17").Value = "  -6.49%  "

$ws.Range("D18").Value = "2.201.06"
$ws.Range("E18").Value = "  -5.78%  "

$ws.Range("D19").Value = "41.559.17"
$ws.Range("E19").Value = "  -5.00%  "

$ws.Range("E20").Value = "  -6.79%  "

$ws.Range("D21").Value = "73.95"
$ws.Range("E21").Value = "  -4.65%  "

$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  -6.47%  "

$ws.Range("D23").Value = "233.78"
$ws.Range("E23").Value = "  -7.06%  "

$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +10.45%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("E26").Value = "  -5.10%  "

$ws.Range("D27").Value = "2.42"
$ws.Range("E27").Value = "  -2.48%  "

$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("D29").Value = "9.78"
$ws.Range("E29").Value = "  -5.40%  "

$ws.Range("D30").Value = "168.25"
$ws.Range("E30").Value = "  -4.07%  "

$ws.Range("D31").Value = "20.35"
$ws.Range("E31").Value = "  -7.97%  "

$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -6.52%  "

$ws.Range("E33").Value = "  -6.53%  "

$ws.Range("E34").Value = "  -3.78%  "

$ws.Range("D35").Value = "5.06"
$ws.Range("E35").Value = "  -4.40%  "

$ws.Range("D36").Value = "'4.60"
$ws.Range("E36").Value = "  -8.23%  "

$ws.Range("D37").Value = "3.91"
$ws.Range("E37").Value = "  +4.83%  "

$ws.Range("D38").Value = "23.13"
$ws.Range("E38").Value = "  +17.77%  "

$ws.Range("E39").Value = "  -4.79%  "

$ws.Range("D40").Value = "0.0271"
$ws.Range("E40").Value = "  +0.37%  "

$ws.Range("D41").Value = "5.84"
$ws.Range("E41").Value = "  -8.62%  "

$ws.Range("D42").Value = "65.23"
$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").Value = "8.87"
$ws.Range("E43").Value = "  -1.88%  "

$ws.Range("D44").Value = "4.81"
$ws.Range("E44").Value = "  -11.01%  "

$ws.Range("D45").Value = "0.192"
$ws.Range("E45").Value = "  -2.25%  "

$ws.Range("E46").Value = "  -5.85%  "

$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").Value = "'4.50"
$ws.Range("E48").Value = "  +4.99%  "

$ws.Range("D49").Value = "10.21"
$ws.Range("E49").Value = "  +8.43%  "

$ws.Range("E50").Value = "  -4.02%  "

$ws.Range("E51").Value = "  +9.69%  "
